# The exported data used to include a header row ("Key" / "Value") in A1:B1,
# pushing the actual key/value pairs down to start at A2. The Excel output
# generation was changed to use header=False, so the header row is removed
# and all data now starts at row 1 (A1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(1).Delete()
